# Auto-generated edit script: updates market-price derived columns (H-N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the scheduled-runner price refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 390
$ws.Range("I12").Value = 390
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 390
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -220
$ws.Range("N12").ClearContents()
$ws.Range("H28").Value = 1257.75
$ws.Range("I28").Value = 1224.8572
$ws.Range("K28").Value = 1224.8572
$ws.Range("M28").Value = -739.8571999999999
$ws.Range("H49").Value = 5272.5
$ws.Range("I49").Value = 545
$ws.Range("J49").Value = 10000
$ws.Range("K49").Value = 1635
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = -1499
$ws.Range("N49").Value = -30272
$ws.Range("H62").Value = 4730
$ws.Range("J62").Value = 3996.6667
$ws.Range("L62").Value = 3996.6667
$ws.Range("N62").Value = -5244.6667
$ws.Range("H65").Value = 4730
$ws.Range("J65").Value = 3996.6667
$ws.Range("L65").Value = 19983.3335
$ws.Range("N65").Value = -26223.3335
$ws.Range("H98").Value = 2218.3635
$ws.Range("J98").Value = 2966
$ws.Range("L98").Value = 2966
$ws.Range("N98").Value = -5962
$ws.Range("H103").Value = 1466.3334
$ws.Range("I103").Value = 1200
$ws.Range("K103").Value = 3600
$ws.Range("M103").Value = -3014
$ws.Range("H113").Value = 8262.444
$ws.Range("I113").Value = 6584.231
$ws.Range("J113").Value = 9820.786
$ws.Range("K113").Value = 6584.231
$ws.Range("L113").Value = 9820.786
$ws.Range("M113").Value = -3330.231
$ws.Range("N113").Value = -16328.786
$ws.Range("H122").Value = 2218.3635
$ws.Range("J122").Value = 2966
$ws.Range("L122").Value = 8898
$ws.Range("N122").Value = -13798
$ws.Range("H135").Value = 23812200
$ws.Range("I135").Value = 29414306
$ws.Range("K135").Value = 264728754
$ws.Range("M135").Value = -264726219

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2693.077
$ws.Range("I45").Value = 1718.8334
$ws.Range("J45").Value = 3528.1428
$ws.Range("K45").Value = 1718.8334
$ws.Range("L45").Value = 3528.1428
$ws.Range("M45").Value = -1341.8334
$ws.Range("N45").Value = -4282.1428
$ws.Range("H74").Value = 1458.6046
$ws.Range("I74").Value = 1474.4359
$ws.Range("K74").Value = 1474.4359
$ws.Range("M74").Value = -600.4358999999999
$ws.Range("H77").Value = 1458.6046
$ws.Range("I77").Value = 1474.4359
$ws.Range("K77").Value = 7372.1795
$ws.Range("M77").Value = -3004.1795
$ws.Range("H97").Value = 918.92
$ws.Range("J97").Value = 909.375
$ws.Range("L97").Value = 909.375
$ws.Range("N97").Value = -1901.375
$ws.Range("H110").Value = 4170.263
$ws.Range("I110").Value = 4124.1665
$ws.Range("K110").Value = 4124.1665
$ws.Range("M110").Value = -2079.1665
$ws.Range("H122").Value = 4257.4517
$ws.Range("I122").Value = 3697.5881
$ws.Range("J122").Value = 4937.2856
$ws.Range("K122").Value = 11092.7643
$ws.Range("L122").Value = 14811.8568
$ws.Range("M122").Value = -8642.764299999999
$ws.Range("N122").Value = -19711.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19823.666
$ws.Range("I26").Value = 19823.666
$ws.Range("K26").Value = 19823.666
$ws.Range("M26").Value = -19531.666
$ws.Range("H94").Value = 2100.6829
$ws.Range("J94").Value = 2548.6667
$ws.Range("L94").Value = 2548.6667
$ws.Range("N94").Value = -3450.6667
$ws.Range("H107").Value = 9981.632
$ws.Range("I107").Value = 6581.357
$ws.Range("K107").Value = 6581.357
$ws.Range("M107").Value = -4661.357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7102.4707
$ws.Range("I16").Value = 7767.4287
$ws.Range("J16").Value = 3999.3333
$ws.Range("K16").Value = 7767.4287
$ws.Range("L16").Value = 3999.3333
$ws.Range("M16").Value = -7480.4287
$ws.Range("N16").Value = -4573.3333
$ws.Range("H31").Value = 1591.5358
$ws.Range("I31").Value = 1591.5358
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1591.5358
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1296.5358
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1591.5358
$ws.Range("I34").Value = 1591.5358
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1591.5358
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1389.5358
$ws.Range("N34").ClearContents()
$ws.Range("H62").Value = 4147.067
$ws.Range("J62").Value = 3949.8
$ws.Range("L62").Value = 3949.8
$ws.Range("N62").Value = -5197.8
$ws.Range("H65").Value = 4147.067
$ws.Range("J65").Value = 3949.8
$ws.Range("L65").Value = 19749
$ws.Range("N65").Value = -25989
$ws.Range("H94").Value = 9780.083000000001
$ws.Range("J94").Value = 1778.25
$ws.Range("L94").Value = 1778.25
$ws.Range("N94").Value = -2680.25
$ws.Range("H109").Value = 20047.334
$ws.Range("J109").Value = 20047.334
$ws.Range("L109").Value = 20047.334
$ws.Range("N109").Value = -22127.334
$ws.Range("H113").Value = 7102.4707
$ws.Range("I113").Value = 7767.4287
$ws.Range("J113").Value = 3999.3333
$ws.Range("K113").Value = 7767.4287
$ws.Range("L113").Value = 3999.3333
$ws.Range("M113").Value = -5597.4287
$ws.Range("N113").Value = -8339.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1034.4736
$ws.Range("I5").Value = 927.6923
$ws.Range("K5").Value = 2783.0769
$ws.Range("M5").Value = -2671.0769
$ws.Range("H47").Value = 819.61536
$ws.Range("I47").Value = 64.2
$ws.Range("K47").Value = 192.6
$ws.Range("M47").Value = 238.4
$ws.Range("H95").Value = 6633.3335
$ws.Range("J95").Value = 6633.3335
$ws.Range("L95").Value = 19900.0005
$ws.Range("N95").Value = -24018.0005
$ws.Range("H135").Value = 1034.4736
$ws.Range("I135").Value = 927.6923
$ws.Range("K135").Value = 8349.2307
$ws.Range("M135").Value = -5814.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 500
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H12").Value = 5000
$ws.Range("I12").Value = 5000
$ws.Range("K12").Value = 5000
$ws.Range("M12").Value = -4860
$ws.Range("H122").Value = 12400
$ws.Range("I122").Value = 25000
$ws.Range("J122").Value = 9250
$ws.Range("K122").Value = 75000
$ws.Range("L122").Value = 27750
$ws.Range("M122").Value = -72550
$ws.Range("N122").Value = -32650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 660.1
$ws.Range("I16").Value = 557.4286
$ws.Range("J16").Value = 899.6667
$ws.Range("K16").Value = 557.4286
$ws.Range("L16").Value = 899.6667
$ws.Range("M16").Value = -387.4286
$ws.Range("N16").Value = -1239.6667
$ws.Range("H60").Value = 17000
$ws.Range("J60").Value = 17000
$ws.Range("L60").Value = 17000
$ws.Range("N60").Value = -18018
$ws.Range("H61").Value = 1490.6875
$ws.Range("I61").Value = 1490.6875
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1490.6875
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1288.6875
$ws.Range("N61").ClearContents()
$ws.Range("H82").Value = 1721.45
$ws.Range("I82").Value = 1829.5714
$ws.Range("K82").Value = 1829.5714
$ws.Range("M82").Value = -1468.5714
$ws.Range("H85").Value = 1721.45
$ws.Range("I85").Value = 1829.5714
$ws.Range("K85").Value = 1829.5714
$ws.Range("M85").Value = -581.5714
$ws.Range("H100").Value = 106202.09
$ws.Range("I100").Value = 187659.5
$ws.Range("J100").Value = 8453.200000000001
$ws.Range("K100").Value = 187659.5
$ws.Range("L100").Value = 8453.200000000001
$ws.Range("M100").Value = -187118.5
$ws.Range("N100").Value = -9535.200000000001
$ws.Range("H113").Value = 1490.6875
$ws.Range("I113").Value = 1490.6875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1490.6875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 679.3125
$ws.Range("N113").ClearContents()
$ws.Range("H117").Value = 73000
$ws.Range("J117").Value = 73000
$ws.Range("L117").Value = 73000
$ws.Range("N117").Value = -82178
$ws.Range("H122").Value = 4928.933
$ws.Range("I122").Value = 3570.5557
$ws.Range("J122").Value = 6966.5
$ws.Range("K122").Value = 10711.6671
$ws.Range("L122").Value = 20899.5
$ws.Range("M122").Value = -8261.667099999999
$ws.Range("N122").Value = -25799.5
$ws.Range("H132").Value = 2458.476
$ws.Range("I132").Value = 2162.4375
$ws.Range("J132").Value = 3405.8
$ws.Range("K132").Value = 6487.3125
$ws.Range("L132").Value = 10217.4
$ws.Range("M132").Value = -3957.3125
$ws.Range("N132").Value = -15277.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4134.0835
$ws.Range("J96").Value = 4093.7144
$ws.Range("L96").Value = 4093.7144
$ws.Range("N96").Value = -6839.7144
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H136").Value = 1602.0588
$ws.Range("I136").Value = 1345.4286
$ws.Range("J136").Value = 2799.6667
$ws.Range("K136").Value = 4036.2858
$ws.Range("L136").Value = 8399.000100000001
$ws.Range("M136").Value = -1486.2858
$ws.Range("N136").Value = -13499.0001
